$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: add missing PriceChange (X5) and UpDown (Y5) values ---
$ws.Range("X4:Y4").Copy()
$ws.Range("X5:Y5").PasteSpecial(-4122)
$ws.Range("X5").Value = -0.59999799999999937
$ws.Range("Y5").Value = "Down"

# --- Row 6: new data row copied (formats) from row 5, then values set ---
$ws.Range("A5:W5").Copy()
$ws.Range("A6:W6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42647.883217592593
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 20456
$ws.Range("F6").Value = 1059
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = 48
$ws.Range("I6").Value = 77
$ws.Range("J6").Value = 21
$ws.Range("K6").Value = 22706
$ws.Range("L6").Value = 200
$ws.Range("M6").Value = 204
$ws.Range("N6").Value = 46
$ws.Range("O6").Value = 13
$ws.Range("P6").Value = "Named"
$ws.Range("Q6").Value = 53.235658945584888
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = -0.0862
$ws.Range("T6").Value = -0.0166
$ws.Range("U6").Value = 6.69
$ws.Range("V6").Value = 1.88
$ws.Range("W6").Value = 0
